$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Build the Excel Table ("Table1") over A1:U75 -------------------------
# A ListObject is created on a small, unformatted scratch range first so
# that the engine does not capture the existing bold/shaded header style
# (A1:U1 already carries cell style index 1) as a table "headerRowDxfId".
# The table is then resized onto the real data range, which keeps the
# worksheet's existing cell formatting/styles completely untouched.
$ws.Range("W1").Value = "tmpHeader"
$ws.Range("W2").Value = "tmpData"
$tbl = $ws.ListObjects.Add(1, $ws.Range("W1:W2"), [System.Reflection.Missing]::Value, 1)
$tbl.Name = "Table1"
$tbl.Resize($ws.Range("A1:U75"))
$ws.Range("W1:W2").Clear()

# --- Rename the "_old" / "_new" header suffixes to "_FV2210" / "_FV2304" --
# Writing straight into the header-row cells both updates the visible
# worksheet text and (because those cells are the table's header row)
# resynchronizes the table's ListColumn names to match.
$ws.Range("A1").Value = "Segmentname_FV2210"
$ws.Range("B1").Value = "Segmentgruppe_FV2210"
$ws.Range("C1").Value = "Segment_FV2210"
$ws.Range("D1").Value = "Datenelement_FV2210"
$ws.Range("E1").Value = "Segment ID_FV2210"
$ws.Range("F1").Value = "Code_FV2210"
$ws.Range("G1").Value = "Qualifier_FV2210"
$ws.Range("H1").Value = "Beschreibung_FV2210"
$ws.Range("I1").Value = "Bedingungsausdruck_FV2210"
$ws.Range("J1").Value = "Bedingung_FV2210"
$ws.Range("K1").Value = "diff"
$ws.Range("L1").Value = "Segmentname_FV2304"
$ws.Range("M1").Value = "Segmentgruppe_FV2304"
$ws.Range("N1").Value = "Segment_FV2304"
$ws.Range("O1").Value = "Datenelement_FV2304"
$ws.Range("P1").Value = "Segment ID_FV2304"
$ws.Range("Q1").Value = "Code_FV2304"
$ws.Range("R1").Value = "Qualifier_FV2304"
$ws.Range("S1").Value = "Beschreibung_FV2304"
$ws.Range("T1").Value = "Bedingungsausdruck_FV2304"
$ws.Range("U1").Value = "Bedingung_FV2304"

# --- Freeze the header row (pane split after row 1) -----------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
